$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 438019870
$ws.Range("J3").Value = 438019870
$ws.Range("L3").Value = 438019870
$ws.Range("N3").Value = -438020098
$ws.Range("H6").Value = 1362.2778
$ws.Range("I6").Value = 205.33333
$ws.Range("J6").Value = 1940.75
$ws.Range("K6").Value = 615.99999
$ws.Range("L6").Value = 5822.25
$ws.Range("M6").Value = -503.99999
$ws.Range("N6").Value = -6046.25
$ws.Range("H12").Value = 200.5
$ws.Range("I12").Value = 200.5
$ws.Range("K12").Value = 200.5
$ws.Range("M12").Value = -30.5
$ws.Range("H48").Value = 4948.4365
$ws.Range("I48").Value = 1053
$ws.Range("J48").Value = 5020.574
$ws.Range("K48").Value = 3159
$ws.Range("L48").Value = 15061.722
$ws.Range("M48").Value = -2867
$ws.Range("N48").Value = -15645.722
$ws.Range("H56").Value = 4948.4365
$ws.Range("I56").Value = 1053
$ws.Range("J56").Value = 5020.574
$ws.Range("K56").Value = 3159
$ws.Range("L56").Value = 15061.722
$ws.Range("M56").Value = -2625
$ws.Range("N56").Value = -16129.722
$ws.Range("H100").Value = 16766476
$ws.Range("I100").Value = 31375706
$ws.Range("K100").Value = 31375706
$ws.Range("M100").Value = -31375165
$ws.Range("H102").Value = 438019870
$ws.Range("J102").Value = 438019870
$ws.Range("L102").Value = 438019870
$ws.Range("N102").Value = -438026360
$ws.Range("H107").Value = 16217.875
$ws.Range("I107").Value = 20340.5
$ws.Range("K107").Value = 20340.5
$ws.Range("M107").Value = -18420.5
$ws.Range("H125").Value = 4589.3076
$ws.Range("I125").Value = 5340.75
$ws.Range("K125").Value = 48066.75
$ws.Range("M125").Value = -45606.75
$ws.Range("H129").Value = 1102.7059
$ws.Range("I129").Value = 1049
$ws.Range("J129").Value = 1277.25
$ws.Range("K129").Value = 3147
$ws.Range("L129").Value = 3831.75
$ws.Range("M129").Value = 1853
$ws.Range("N129").Value = -13831.75
$ws.Range("H137").Value = 9865.677
$ws.Range("I137").Value = 13860.65
$ws.Range("J137").Value = 4158.5713
$ws.Range("K137").Value = 41581.95
$ws.Range("L137").Value = 12475.7139
$ws.Range("M137").Value = -39031.95
$ws.Range("N137").Value = -17575.7139
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4633.4907
$ws.Range("I32").Value = 4684.51
$ws.Range("K32").Value = 4684.51
$ws.Range("M32").Value = -4397.51
$ws.Range("H36").Value = 13333.167
$ws.Range("I36").Value = 11999.8
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 11999.8
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -11653.8
$ws.Range("N36").Value = -20692
$ws.Range("H94").Value = 188605540
$ws.Range("J94").Value = 188605540
$ws.Range("L94").Value = 188605540
$ws.Range("N94").Value = -188607342
$ws.Range("H110").Value = 2217
$ws.Range("I110").Value = 611.6667
$ws.Range("K110").Value = 611.6667
$ws.Range("M110").Value = 1433.3333
$ws.Range("H132").Value = 6738.8203
$ws.Range("I132").Value = 7075.483
$ws.Range("K132").Value = 21226.449
$ws.Range("M132").Value = -18696.449
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 13539.121
$ws.Range("I99").Value = 18438.715
$ws.Range("J99").Value = 4964.8335
$ws.Range("K99").Value = 18438.715
$ws.Range("L99").Value = 4964.8335
$ws.Range("M99").Value = -16940.715
$ws.Range("N99").Value = -7960.8335
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1104.5555
$ws.Range("I22").Value = 823
$ws.Range("K22").Value = 823
$ws.Range("M22").Value = -473
$ws.Range("H31").Value = 13408.143
$ws.Range("I31").Value = 32458.25
$ws.Range("K31").Value = 32458.25
$ws.Range("M31").Value = -32163.25
$ws.Range("H34").Value = 13408.143
$ws.Range("I34").Value = 32458.25
$ws.Range("K34").Value = 32458.25
$ws.Range("M34").Value = -32256.25
$ws.Range("H99").Value = 135764.81
$ws.Range("I99").Value = 297450.28
$ws.Range("K99").Value = 297450.28
$ws.Range("M99").Value = -295952.28
$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H126").Value = 135764.81
$ws.Range("I126").Value = 297450.28
$ws.Range("K126").Value = 892350.8400000001
$ws.Range("M126").Value = -889880.8400000001
$ws.Range("H132").Value = 17122.037
$ws.Range("I132").Value = 2373.4443
$ws.Range("K132").Value = 7120.3329
$ws.Range("M132").Value = -4590.3329
$ws.Range("H134").Value = 3980.7273
$ws.Range("I134").Value = 4953.857
$ws.Range("K134").Value = 14861.571
$ws.Range("M134").Value = -12326.571
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19989160
$ws.Range("I4").Value = 3053813.2
$ws.Range("K4").Value = 9161439.600000001
$ws.Range("M4").Value = -9161327.600000001
$ws.Range("H34").Value = 2990.5
$ws.Range("I34").Value = 2573.25
$ws.Range("J34").Value = 3825
$ws.Range("K34").Value = 7719.75
$ws.Range("L34").Value = 11475
$ws.Range("M34").Value = -7635.75
$ws.Range("N34").Value = -11643
$ws.Range("H55").Value = 4266.722
$ws.Range("J55").Value = 4893.1333
$ws.Range("L55").Value = 14679.3999
$ws.Range("N55").Value = -15033.3999
$ws.Range("H80").Value = 49464.61
$ws.Range("J80").Value = 72377
$ws.Range("L80").Value = 217131
$ws.Range("N80").Value = -219003
$ws.Range("H83").Value = 49464.61
$ws.Range("J83").Value = 72377
$ws.Range("L83").Value = 651393
$ws.Range("N83").Value = -660753
$ws.Range("H121").Value = 1113179.5
$ws.Range("I121").Value = 622.1667
$ws.Range("J121").Value = 1669458.1
$ws.Range("K121").Value = 1866.5001
$ws.Range("L121").Value = 5008374.300000001
$ws.Range("M121").Value = -556.5001
$ws.Range("N121").Value = -5010994.300000001
$ws.Range("H131").Value = 1180.9231
$ws.Range("J131").Value = 2499.5
$ws.Range("L131").Value = 7498.5
$ws.Range("N131").Value = -17578.5
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15068.286
$ws.Range("I80").Value = 25522
$ws.Range("J80").Value = 1130
$ws.Range("K80").Value = 25522
$ws.Range("L80").Value = 1130
$ws.Range("M80").Value = -24524
$ws.Range("N80").Value = -3126
$ws.Range("H83").Value = 15068.286
$ws.Range("I83").Value = 25522
$ws.Range("J83").Value = 1130
$ws.Range("K83").Value = 127610
$ws.Range("L83").Value = 5650
$ws.Range("M83").Value = -122618
$ws.Range("N83").Value = -15634
$ws.Range("H102").Value = 5595.515
$ws.Range("I102").Value = 7628.2856
$ws.Range("K102").Value = 7628.2856
$ws.Range("M102").Value = -6006.2856
$ws.Range("H104").Value = 67100
$ws.Range("J104").Value = 67100
$ws.Range("L104").Value = 67100
$ws.Range("N104").Value = -74088
$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140
$ws.Range("H137").Value = 129946.75
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 129946.75
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 129946.75
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -140146.75
$ws.Range("H138").Value = 119749
$ws.Range("J138").Value = 119749
$ws.Range("L138").Value = 119749
$ws.Range("N138").Value = -130029
$ws.Range("H139").Value = 47895.75
$ws.Range("J139").Value = 47895.75
$ws.Range("L139").Value = 47895.75
$ws.Range("N139").Value = -58175.75
$ws.Range("H140").Value = 77441.60000000001
$ws.Range("J140").Value = 86624.75
$ws.Range("L140").Value = 86624.75
$ws.Range("N140").Value = -96984.75
$ws.Range("H141").Value = 62808.332
$ws.Range("J141").Value = 62785.375
$ws.Range("L141").Value = 62785.375
$ws.Range("N141").Value = -73145.375
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5695.8423
$ws.Range("I122").Value = 4863.154
$ws.Range("K122").Value = 14589.462
$ws.Range("M122").Value = -12139.462
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 12390.333
$ws.Range("J69").Value = 12390.333
$ws.Range("L69").Value = 12390.333
$ws.Range("N69").Value = -13888.333
$ws.Range("H72").Value = 12390.333
$ws.Range("J72").Value = 12390.333
$ws.Range("L72").Value = 37170.999
$ws.Range("N72").Value = -44658.999
$ws.Range("H104").Value = 13795.5
$ws.Range("J104").Value = 13795.5
$ws.Range("L104").Value = 13795.5
$ws.Range("N104").Value = -20783.5
$ws.Range("H113").Value = 1384.1666
$ws.Range("I113").Value = 645.25
$ws.Range("K113").Value = 1935.75
$ws.Range("M113").Value = 234.25
$ws.Range("H126").Value = 16227.866
$ws.Range("I126").Value = 22104.75
$ws.Range("K126").Value = 66314.25
$ws.Range("M126").Value = -63844.25
